$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the "Meta description: ..." paragraph that currently sits right
#    after the title heading ("Play Aztec Magic Deluxe for Free - Review and
#    Features").
# ---------------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
if ($metaPara.Range.Text -like "Meta description:*") {
    $metaPara.Range.Delete()
}

# ---------------------------------------------------------------------------
# 2) The last paragraph in the document currently holds the image-generation
#    "Prompt: ..." text (italic). Split it into two paragraphs:
#      a) a new bold paragraph with the title text
#         "Play Aztec Magic Deluxe for Free - Review and Features"
#      b) the existing (now last) paragraph, with its text replaced by the
#         meta-description sentence (keeping the italic formatting), i.e.
#         the "Meta description:" prefix is dropped.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)

# Sanity check we are looking at the expected paragraph before editing it.
if ($lastPara.Range.Text -like "Prompt:*") {

    # Insert a brand new empty paragraph right after the second-to-last
    # paragraph (i.e. immediately before the "Prompt:" paragraph).
    $prevPara = $d.Paragraphs.Item($count - 1)
    $prevPara.Range.InsertParagraphAfter() | Out-Null

    # The newly created paragraph is now item($count); make sure it uses the
    # plain "Normal" style (it otherwise would inherit the style of whatever
    # paragraph preceded it).
    $newParaObj = $d.Paragraphs.Item($count)
    $newParaObj.Style = $d.Styles.Item("Normal")

    $newParaRange = $newParaObj.Range
    $titleStart = $newParaRange.Start
    $titleText = "Play Aztec Magic Deluxe for Free - Review and Features"
    $newParaRange.Text = $titleText

    # Bold the title text (without touching the paragraph mark itself, so we
    # don't leave stray bold formatting on the pilcrow).
    $titleRange = $d.Range($titleStart, $titleStart + $titleText.Length)
    $titleRange.Font.Bold = $true

    # The original "Prompt: ..." paragraph is now item($count + 1). Replace
    # its text, dropping the "Prompt: " lead-in, with the meta description
    # sentence, using Find/Replace on that paragraph's range so the existing
    # run (and its italic formatting) is preserved and only the text content
    # changes.
    $descPara = $d.Paragraphs.Item($count + 1)
    $descRange = $descPara.Range
    $oldText = $descRange.Text
    $newText = "Discover the Aztec-themed slot game Aztec Magic Deluxe and play it for free. Read our review of the features and gameplay of this beautifully designed game."

    $find = $descRange.Find
    $find.ClearFormatting()
    $find.Replacement.ClearFormatting()
    $find.Execute($oldText.TrimEnd(), $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null
}
